# Actualizacion de la planilla de control: se completa la consigna de
# subida del sitio (fila 10): pasa de "En progreso" a "SI" y el comentario
# se reemplaza por el enlace publicado del sitio (con hipervinculo).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B10: "En progreso" -> "SI" (misma presentacion que el resto de la
#          columna "Cumplido": relleno verde, texto centrado) ---
$ws.Range("B10").Value = "SI"
$ws.Range("B10").Interior.Color = 5296274
$ws.Range("B10").HorizontalAlignment = -4108
$ws.Range("B10").VerticalAlignment = -4108

# --- C10: comentario de texto libre -> enlace al sitio publicado ---
$ws.Range("C10").Value = "https://sportspro.netlify.app/"
$ws.Hyperlinks.Add($ws.Range("C10"), "https://sportspro.netlify.app/")

# La vista queda desplazada mostrando la fila recien actualizada, con
# B10 como celda activa.
[void]$ws.Range("B10").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1

$wb.Save()
